# ============================================================================
# Edit script: applies the changes described by the commit
#   "📝 PPT 2 to 4, index 수정"
# against "02.자바스크립트 - 함수.pptx"
#
# Changes:
#  1. Auto date fields (footer placeholders in the slide master + the 11
#     slide layouts, type "datetime1") bump from 2025-04-22 -> 2025-04-23.
#  2. Slide 15, "TextBox 7" (console-output box): the arrow-function demo
#     output is corrected from `11 ＇바로핑＇ / 9 ＇라라핑'`
#     to `11 '' / 9 ''` (this.name is undefined inside the arrow function).
#  3. Slide 16, "Text Box 3": spelling fix 메소드의 -> 메서드의.
#  4. Slide 23, "TextBox 8": code comment fixed from `// 30` to `// 70`
#     (index correction referenced in the commit message).
# ============================================================================

$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# 1) Date placeholders: slide master + all 11 slide layouts
# ----------------------------------------------------------------------
$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = "2025-04-23"

$layoutDateShape = @{1=3; 2=3; 3=3; 4=4; 5=6; 6=2; 7=1; 8=4; 9=4; 10=3; 11=3}
for ($li = 1; $li -le 11; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $idx = $layoutDateShape[$li]
    $layout.Shapes.Item($idx).TextFrame.TextRange.Text = "2025-04-23"
}

# NotesMaster footer date ("4/22/2025" -> "4/23/2025") lives on the
# notesMaster's "datetimeFigureOut" field.
$p.NotesMaster.Shapes.Item(2).TextFrame.TextRange.Text = "4/23/2025"

# ----------------------------------------------------------------------
# 2) Slide 15 - "TextBox 7" console output text fix
# ----------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$tb7 = $s15.Shapes.Item(5)
$tb7tr = $tb7.TextFrame.TextRange
# Paragraph 1: "11 ＇바로핑＇" (8 chars) -> "11 ''"
$tb7tr.Paragraphs(1, 1).Characters(1, 8).Text = "11 ''"
# Paragraph 2: "9 ＇라라핑'" - keep the leading "9", replace the rest
# (chars 2-7 relative to paragraph start: " ＇라라핑'") with " ''"
$tb7tr2 = $tb7.TextFrame.TextRange
$tb7tr2.Paragraphs(2, 1).Characters(2, 6).Text = " ''"

# ----------------------------------------------------------------------
# 3) Slide 16 - "Text Box 3" spelling fix (메소드의 -> 메서드의)
# ----------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$tb3 = $s16.Shapes.Item(4)
$tb3tr = $tb3.TextFrame.TextRange
$tb3tr.Paragraphs(4, 1).Characters(24, 4).Text = "메서드의"

# ----------------------------------------------------------------------
# 4) Slide 23 - "TextBox 8" code comment fix (// 30 -> // 70)
# ----------------------------------------------------------------------
$s23 = $p.Slides.Item(23)
$tb8 = $s23.Shapes.Item(5)
$tb8tr = $tb8.TextFrame.TextRange
$tb8tr.Paragraphs(6, 1).Runs(7, 1).Text = "// 70"
